# Fruta / hortaliza, semanal
# Insert two new weekly price rows into the "Vega Modelo de Temuco - Damasco" sheet.
# Row 24 (new) and row 40 (new, after the first insert shifts everything down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 24 - shifts current rows 24.. down by one.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24.
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 44554
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100103
$ws.Range("H24").Value = "Frutos de hueso (carozo)"
$ws.Range("I24").Value = 100103003
$ws.Range("J24").Value = "Damasco"
$ws.Range("K24").Value = "Dina"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 35
$ws.Range("N24").Value = 20000
$ws.Range("O24").Value = 20000
$ws.Range("P24").Value = 20000
$ws.Range("Q24").Value = '$/caja 18 kilos'
$ws.Range("R24").Value = "Provincia de Quillota"
$ws.Range("S24").Value = 1111
$ws.Range("T24").Value = 18

# Insert a second blank row at position 40 - shifts current rows 40.. down by one.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40.
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "Vega Modelo de Temuco"
$ws.Range("C40").Value = "La Araucanía"
$ws.Range("D40").Value = 44553
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100103
$ws.Range("H40").Value = "Frutos de hueso (carozo)"
$ws.Range("I40").Value = 100103003
$ws.Range("J40").Value = "Damasco"
$ws.Range("K40").Value = "Dina"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 65
$ws.Range("N40").Value = 20000
$ws.Range("O40").Value = 20000
$ws.Range("P40").Value = 20000
$ws.Range("Q40").Value = '$/caja 18 kilos'
$ws.Range("R40").Value = "Provincia de Quillota"
$ws.Range("S40").Value = 1111
$ws.Range("T40").Value = 18
